$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "262k02"
$ws.Range("B5").Value = "Al"
$ws.Range("C5").Value = "Dabri"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0

$ws.Range("H5").Select()
